$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 10:48:02"
$wsZhCn.Range("H2").Value = "2016-03-12 10:48:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 10:48:05"
$wsDeDe.Range("H2").Value = "2016-03-12 10:48:25"
